$d = $word.ActiveDocument

# The text "<fr>...</fr>" wrapping is used many times throughout this
# document, so a document-wide Find/Replace of "</fr>" (or similar)
# would risk touching the wrong occurrence. Scope the edit to the one
# paragraph that contains the phrase we care about.
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Eau de vye*") {
        $para = $p
        break
    }
}

# 1) Drop the closing "</fr>" tag entirely -- its whole run disappears.
$rng = $para.Range
$null = $rng.Find.Execute("</fr>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 2)

# 2) Collapse the opening tag run from "<m><fr>" down to just "<m>".
$rng = $para.Range
$null = $rng.Find.Execute("<m><fr>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<m>", 2)

# 3) Fix the misspelled word in its own run: "Eau de vye" -> "Eau-de-vie".
$rng = $para.Range
$null = $rng.Find.Execute("Eau de vye", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Eau-de-vie", 2)
